# Hortaliza, Femacal de La Calera - Haba
# A new weekly price observation was inserted as row 40, pushing the
# existing rows 40-112 down to 41-113 (dimension grows from A1:R112 to
# A1:R113).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 40 - shifts rows 40:112 -> 41:113
$ws.Rows("40:40").Insert()

# Populate the new row 40 with the new observation
$ws.Range("A40").Value = 3
$ws.Range("B40").Value = "Femacal de La Calera"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44544
$ws.Range("E40").Value = 5
$ws.Range("F40").Value = 100112026
$ws.Range("G40").Value = "Haba"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 85
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 8500
$ws.Range("M40").Value = 8235
$ws.Range("N40").Value = "`$/malla 25 kilos"
$ws.Range("O40").Value = "Provincia de Limarí"
$ws.Range("P40").Value = 329
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"
